$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-11 (Original category, Complementary category, frequency)
$data = @(
    @("Home Decor", "Home Decor", 15848),
    @("Kitchen & Dining", "Kitchen & Dining", 10976),
    @("Seasonal & Holidays", "Seasonal & Holidays", 7946),
    @("Fashion & Accessories", "Fashion & Accessories", 6008),
    @("Stationery & Office", "Stationery & Office", 5858),
    @("Home Decor", "Kitchen & Dining", 5347),
    @("Kitchen & Dining", "Home Decor", 5347),
    @("Home Decor", "Seasonal & Holidays", 4392),
    @("Seasonal & Holidays", "Home Decor", 4392),
    @("Fashion & Accessories", "Home Decor", 3057)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
